# [GQA] Correção dos relatórios de NC.
$wb = $excel.ActiveWorkbook

# --- Sheet "Processo" ---
$wsProcesso = $wb.Worksheets.Item("Processo")
$wsProcesso.Range("B20").Value2 = "Fora da data de avaliação"

# --- Sheet "Análise" ---
$wsAnalise = $wb.Worksheets.Item("Análise")
$wsAnalise.Range("B19").Value2 = "Sim"
$wsAnalise.Range("C19").Value2 = ""
$wsAnalise.Range("D19").Value2 = ""

# --- Sheet "Codificação" (selection only, text content index renumbers itself) ---
$wsCodificacao = $wb.Worksheets.Item("Codificação")

# --- Sheet "Teste" ---
$wsTeste = $wb.Worksheets.Item("Teste")

$b12 = $wsTeste.Range("B12")
$b12.Value2 = "Fora da data de avaliação"
$b12.Font.Name = "Arial"
$b12.Font.Size = 12
$b12.Font.Bold = $false
$b12.Font.ColorIndex = 1
$b12.Interior.Pattern = -4142
$b12.HorizontalAlignment = -4108
$b12.WrapText = $true
$b12.Borders.LineStyle = 1
$b12.Borders.Weight = 2

$b13 = $wsTeste.Range("B13")
$b13.Value2 = "Fora da data de avaliação"
$b13.Font.Name = "Arial"
$b13.Font.Size = 12
$b13.Font.Bold = $false
$b13.Font.ColorIndex = 1
$b13.Interior.Pattern = -4142
$b13.HorizontalAlignment = -4108
$b13.WrapText = $true
$b13.Borders.LineStyle = 1
$b13.Borders.Weight = 2

$wsTeste.Range("C13").Value2 = ""
$wsTeste.Range("D13").Value2 = ""

# --- Sheet "Legenda" ---
$wsLegenda = $wb.Worksheets.Item("Legenda")

# --- Selections per sheet ---
$wsProcesso.Activate()
$wsProcesso.Range("D37").Select() | Out-Null

$wsAnalise.Activate()
$wsAnalise.Range("D19").Select() | Out-Null

$wsProjeto = $wb.Worksheets.Item("Projeto")
$wsProjeto.Activate()
$wsProjeto.Range("A27").Select() | Out-Null

$wsCodificacao.Activate()
$wsCodificacao.Range("A36").Select() | Out-Null

$wsLegenda.Activate()
$wsLegenda.Range("B3").Select() | Out-Null

# Teste becomes the active/selected sheet (activeTab=4), with new selection B13
$wsTeste.Activate()
$wsTeste.Range("B13").Select() | Out-Null
